$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header changes
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Row 2
$ws.Range("D2").Value = "wRLS"
$ws.Range("E2").Value = 163.0355370277466
$ws.Range("F2").Value = 0.3553169473669019
$ws.Range("G2").Value = 127.9919714595157

# Row 3
$ws.Range("D3").Value = "wRLS"
$ws.Range("E3").Value = 160.411597115793
$ws.Range("F3").Value = 0.3495983762100451
$ws.Range("G3").Value = 125.9038434470185

# Row 4
$ws.Range("D4").Value = "wRLS"
$ws.Range("E4").Value = 160.6494430559906
$ws.Range("F4").Value = 0.3501167337102277
$ws.Range("G4").Value = 126.115402345334

# Row 5
$ws.Range("D5").Value = "wRLS"
$ws.Range("E5").Value = 164.2488387549586
$ws.Range("F5").Value = 0.3579611970428152
$ws.Range("G5").Value = 129.9741451310439

# Row 6
$ws.Range("D6").Value = "wRLS"
$ws.Range("E6").Value = 164.8130581386089
$ws.Range("F6").Value = 0.3591908474165844
$ws.Range("G6").Value = 130.4966806033764

# Row 7
$ws.Range("D7").Value = "wRLS"
$ws.Range("E7").Value = 164.9010690741356
$ws.Range("F7").Value = 0.3593826569908429
$ws.Range("G7").Value = 130.5834700776883

# Row 8
$ws.Range("D8").Value = "wRLS"
$ws.Range("E8").Value = 163.2808532193095
$ws.Range("F8").Value = 0.3558515854091036
$ws.Range("G8").Value = 129.0280619757758
